$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 565
$ws.Range("I2").Value = 1540
$ws.Range("J2").Value = 6369
$ws.Range("K2").Value = 41
$ws.Range("L2").Value = 1746
$ws.Range("M2").Value = 114
$ws.Range("N2").Value = 1078
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 28
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 70
$ws.Range("S2").Value = 708
$ws.Range("T2").Value = 1068
$ws.Range("U2").Value = 85
$ws.Range("V2").Value = 9746
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 9737
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 166
$ws.Range("AA2").Value = 80

$wb.Save()
